# Update generated numbers on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

function Update-Sheet {
    param($ws, $f16Value)

    $ws.Cells.Item(2, 6).Value = 11598
    $ws.Cells.Item(3, 6).Value = 11140
    $ws.Cells.Item(6, 6).Value = 1007
    $ws.Cells.Item(11, 6).Value = 10677
    $ws.Cells.Item(12, 6).Value = 4127

    # Row 13: F went from 1 -> 11, and G switched from text "不可售" to number 30.
    $ws.Cells.Item(13, 6).Value = 11
    $ws.Cells.Item(13, 7).Value = 30

    # Row 16: F value differs slightly between the two sheets.
    $ws.Cells.Item(16, 6).Value = $f16Value

    $ws.Cells.Item(17, 6).Value = 45
    $ws.Cells.Item(19, 6).Value = 433
    $ws.Cells.Item(20, 6).Value = 11116
    $ws.Cells.Item(21, 6).Value = 10883
}

$wsExhibition = $wb.Worksheets.Item("展览")
Update-Sheet $wsExhibition 120

$wsAll = $wb.Worksheets.Item("全部类型")
Update-Sheet $wsAll 121

# Row 15 ("F15") only changes on the "全部类型" sheet.
$wsAll.Cells.Item(15, 6).Value = 2461
